# edit.ps1
# Applies the commit "edit projects, temperature defn, KM's role":
#   1. ColumnHeadersNcp!B7 - fix the temperature row's definition text
#      (it previously, incorrectly, described "salinity").
#   2. Personnel - add a new row for Kate Morkeski (metadata Provider,
#      OCE-2322676 project).

$wb = $excel.ActiveWorkbook

# --- 1. Fix temperature definition on the ColumnHeadersNcp sheet ---------
$ncp = $wb.Worksheets.Item("ColumnHeadersNcp")
$ncp.Range("B7").Value = "Underway thermosalinograph temperature in degrees Celsius. URI http://vocab.nerc.ac.uk/collection/P01/current/TEMPSZ01/"

# --- 2. Add Kate Morkeski to the Personnel sheet --------------------------
$personnel = $wb.Worksheets.Item("Personnel")

$personnel.Range("A10").Value = "Kate"
$personnel.Range("C10").Value = "Morkeski"
$personnel.Range("D10").Value = "Northeast U.S. Shelf LTER"
$personnel.Range("E10").Value = "kmorkeski@whoi.edu"
$personnel.Range("F10").Value = "0000-0002-2903-5851"
$personnel.Range("G10").Value = "metadata Provider"
$personnel.Range("H10").Value = "Northeast U.S. Shelf LTER"
$personnel.Range("I10").Value = "NSF"
$personnel.Range("J10").Value = "OCE-2322676"

# --- 3. Update sheet selections to match the new active cells ------------
[void]$ncp.Range("B7").Select()
[void]$personnel.Range("A10:J10").Select()
